# Add a new "Player Info" worksheet in front of the existing sheets,
# populate it with the player's basic info, and update the
# MATCH_CARD_LINK columns (on the ODI Batting / ODI Bowling sheets) to a
# simpler MATCH_CODE column containing just the numeric match code.

$wb = $excel.ActiveWorkbook

# --- 1. Insert the new "Player Info" sheet before the current first sheet ---
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$playerInfo = $wb.Worksheets.Add($battingSheet)
$playerInfo.Name = "Player Info"

$headers = @("ID", "NAME", "BATTING_HAND", "BOWL_STYLE")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $playerInfo.Cells.Item(1, $i + 1)
    $cell.Value = $headers[$i]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

$playerInfo.Cells.Item(2, 1).Value = "'6790"
$playerInfo.Cells.Item(2, 2).Value = "Agha Salman"
$playerInfo.Cells.Item(2, 3).Value = "Right Handed"
$playerInfo.Cells.Item(2, 4).Value = "Right Arm Off Break"

$playerInfo.Range("A1").Select()

# --- 2. Replace MATCH_CARD_LINK with MATCH_CODE on "ODI Batting" ---
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$battingSheet.Range("D1").Value = "MATCH_CODE"

$battingLinks = @{
    2 = "4634"
    3 = "4638"
    4 = "4641"
    5 = "4686"
    6 = "4688"
    7 = "4690"
}
foreach ($row in $battingLinks.Keys) {
    $battingSheet.Cells.Item($row, 4).Value = "'" + $battingLinks[$row]
}

# --- 3. Replace MATCH_CARD_LINK with MATCH_CODE on "ODI Bowling" ---
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")
$bowlingSheet.Range("B1").Value = "MATCH_CODE"

$bowlingLinks = @{
    2 = "4634"
    3 = "4641"
    4 = "4686"
    5 = "4688"
    6 = "4690"
}
foreach ($row in $bowlingLinks.Keys) {
    $bowlingSheet.Cells.Item($row, 2).Value = "'" + $bowlingLinks[$row]
}
